# Weekly Fruit/Vegetable price update: insert 4 new "Choclo" price records
# (dated 2022-02-03, serial 44595) for "Provincia de Melipilla" at row 430,
# pushing the existing rows 430-438 down to 434-442.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 430 (shifts old rows 430:438 -> 434:442)
$ws.Range("A430:A433").EntireRow.Insert()

# Row 430: Choclero / Primera
$ws.Cells.Item(430, 1).Value = 9
$ws.Cells.Item(430, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(430, 3).Value = "Metropolitana"
$ws.Cells.Item(430, 4).Value = 44595
$ws.Cells.Item(430, 5).Value = 13
$ws.Cells.Item(430, 6).Value = 100112024
$ws.Cells.Item(430, 7).Value = "Choclo"
$ws.Cells.Item(430, 8).Value = "Choclero"
$ws.Cells.Item(430, 9).Value = "Primera"
$ws.Cells.Item(430, 10).Value = 15000
$ws.Cells.Item(430, 11).Value = 130
$ws.Cells.Item(430, 12).Value = 170
$ws.Cells.Item(430, 13).Value = 147
$ws.Cells.Item(430, 14).Value = "$/unidad"
$ws.Cells.Item(430, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(430, 16).Value = 147
$ws.Cells.Item(430, 17).Value = 1
$ws.Cells.Item(430, 18).Value = "Hortaliza"

# Row 431: Choclero / Segunda
$ws.Cells.Item(431, 1).Value = 9
$ws.Cells.Item(431, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(431, 3).Value = "Metropolitana"
$ws.Cells.Item(431, 4).Value = 44595
$ws.Cells.Item(431, 5).Value = 13
$ws.Cells.Item(431, 6).Value = 100112024
$ws.Cells.Item(431, 7).Value = "Choclo"
$ws.Cells.Item(431, 8).Value = "Choclero"
$ws.Cells.Item(431, 9).Value = "Segunda"
$ws.Cells.Item(431, 10).Value = 9000
$ws.Cells.Item(431, 11).Value = 100
$ws.Cells.Item(431, 12).Value = 120
$ws.Cells.Item(431, 13).Value = 113
$ws.Cells.Item(431, 14).Value = "$/unidad"
$ws.Cells.Item(431, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(431, 16).Value = 113
$ws.Cells.Item(431, 17).Value = 1
$ws.Cells.Item(431, 18).Value = "Hortaliza"

# Row 432: Dulce o Americano / Primera
$ws.Cells.Item(432, 1).Value = 9
$ws.Cells.Item(432, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(432, 3).Value = "Metropolitana"
$ws.Cells.Item(432, 4).Value = 44595
$ws.Cells.Item(432, 5).Value = 13
$ws.Cells.Item(432, 6).Value = 100112024
$ws.Cells.Item(432, 7).Value = "Choclo"
$ws.Cells.Item(432, 8).Value = "Dulce o Americano"
$ws.Cells.Item(432, 9).Value = "Primera"
$ws.Cells.Item(432, 10).Value = 23000
$ws.Cells.Item(432, 11).Value = 100
$ws.Cells.Item(432, 12).Value = 120
$ws.Cells.Item(432, 13).Value = 109
$ws.Cells.Item(432, 14).Value = "$/unidad"
$ws.Cells.Item(432, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(432, 16).Value = 109
$ws.Cells.Item(432, 17).Value = 1
$ws.Cells.Item(432, 18).Value = "Hortaliza"

# Row 433: Dulce o Americano / Segunda
$ws.Cells.Item(433, 1).Value = 9
$ws.Cells.Item(433, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(433, 3).Value = "Metropolitana"
$ws.Cells.Item(433, 4).Value = 44595
$ws.Cells.Item(433, 5).Value = 13
$ws.Cells.Item(433, 6).Value = 100112024
$ws.Cells.Item(433, 7).Value = "Choclo"
$ws.Cells.Item(433, 8).Value = "Dulce o Americano"
$ws.Cells.Item(433, 9).Value = "Segunda"
$ws.Cells.Item(433, 10).Value = 8000
$ws.Cells.Item(433, 11).Value = 80
$ws.Cells.Item(433, 12).Value = 80
$ws.Cells.Item(433, 13).Value = 80
$ws.Cells.Item(433, 14).Value = "$/unidad"
$ws.Cells.Item(433, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(433, 16).Value = 80
$ws.Cells.Item(433, 17).Value = 1
$ws.Cells.Item(433, 18).Value = "Hortaliza"
